# Update column G ("K") values on the active worksheet to reflect the
# regenerated strikeout-count (K) data from the re-computed save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 4
    4  = 6
    5  = 6
    6  = 7
    7  = 7
    8  = 8
    9  = 9
    10 = 9
    11 = 5
    12 = 10
    13 = 9
    14 = 6
    15 = 10
    16 = 4
    17 = 5
    18 = 7
    19 = 5
    20 = 1
    21 = 8
    22 = 4
    23 = 4
    24 = 3
    25 = 8
    26 = 4
    27 = 3
    28 = 6
    29 = 5
    30 = 5
    31 = 8
    32 = 4
    33 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
